$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description: ..." paragraph right after the
#    title (Heading1) paragraph, matching the existing document's run
#    pattern (a leading empty run, then a bold run, then a plain run).
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:r/>' +
           '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
           '<w:r><w:t>: Join the fun with Bloxx Arctic, a unique slot game with mirrored grids and an activated block feature. Play for free on all devices.</w:t></w:r>' +
           '</w:p>'

$null = $d.Paragraphs(2).Range.InsertXML($metaXml)

# ---------------------------------------------------------------------
# 2) Remove the now-duplicated "Play Bloxx Arctic for Free - Unique
#    gameplay mechanics" (bold) paragraph near the end of the document.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd() -eq "Play Bloxx Arctic for Free - Unique gameplay mechanics") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 3) Replace the closing italic paragraph's text with the new image
#    prompt, preserving its formatting/run structure.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$r = $lastPara.Range
$textRange = $d.Range($r.Start, $r.End - 1)
$textRange.Text = "Please create a cartoon-style image featuring a happy Maya warrior wearing glasses to fit the game `"Bloxx Arctic`". The image can incorporate elements from the game such as the snow and icicles to create a winter theme. The warrior could be holding a dollar symbol or standing next to the Wild snow girl. The background should feature the starry blue sky of Antarctica. The image should be bright and colorful, with attention-grabbing details that will attract potential players to the game."
